$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.5633756828510093
$ws.Range("J2").Value = 0.5633756828510093
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.309024333333333
$ws.Range("N2").Value = 3.927073
$ws.Range("O2").Value = 0.9436068997599097
$ws.Range("P2").Value = 0.9436068997599097
$ws.Range("Q2").Value = 1.265568652539666
$ws.Range("R2").Value = 11.390117872857
$ws.Range("S2").Value = 0.531605181495163
$ws.Range("T2").Value = 0.531605181495163

$ws.Range("I3").Value = 0.5633756828510093
$ws.Range("J3").Value = 0.5633756828510093
$ws.Range("O3").Value = 0.05639310024009027
$ws.Range("P3").Value = 0.05639310024009028
$ws.Range("S3").Value = 0.03177050135584627
$ws.Range("T3").Value = 0.03177050135584628

$ws.Range("G4").Value = 0.7492863333333334
$ws.Range("H4").Value = 2.247859
$ws.Range("I4").Value = 0.4366243171489907
$ws.Range("J4").Value = 0.4366243171489907
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.309024333333333
$ws.Range("N4").Value = 3.927073
$ws.Range("O4").Value = 0.9436068997599097
$ws.Range("P4").Value = 0.9436068997599097
$ws.Range("Q4").Value = 0.9808340429674444
$ws.Range("R4").Value = 8.827506386706998
$ws.Range("S4").Value = 0.4120017182647467
$ws.Range("T4").Value = 0.4120017182647467

$ws.Range("G5").Value = 0.7492863333333334
$ws.Range("H5").Value = 2.247859
$ws.Range("I5").Value = 0.4366243171489907
$ws.Range("J5").Value = 0.4366243171489907
$ws.Range("O5").Value = 0.05639310024009027
$ws.Range("P5").Value = 0.05639310024009028
$ws.Range("Q5").Value = 0.05861791866722222
$ws.Range("R5").Value = 0.527561268005
$ws.Range("S5").Value = 0.024622598884244
$ws.Range("T5").Value = 0.024622598884244
